# Capitalization of column headings in tables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append " DBFS" to the title in B2
$ws.Range("B2").Value = "LLM modell: gemini-2.0-flash;  DBFS"

# Capitalize the first letter of each column header in row 3
$ws.Range("B3").Value = "Kérdések száma"
$ws.Range("C3").Value = "Embedding  generálásai idő átlaga"
$ws.Range("D3").Value = "Sparse embedding generálási idő átlaga"
$ws.Range("E3").Value = "Kontextus összeállitási idő átlaga"
$ws.Range("F3").Value = "LLM feldolgozási idő átlaga"
$ws.Range("G3").Value = "Teljes feldoldozási idő átlaga"
$ws.Range("H3").Value = "Szemantikus hasonlóság mérékének  (BERTScore F1) átlaga (0-1) között"
$ws.Range("I3").Value = "Top_k darab számának átlag"

# Update the active selection to match the new saved state
$ws.Range("D14").Select()
